# Auto-generated Excel COM-interop edit script
# Applies the Dec 10 2023 cryptos list refresh (GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.817.72"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.356.84"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'239.97"
$ws.Range("D6").Value = "'0.669"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").Value = "'73.33"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").Value = "'60.84"
$ws.Range("E11").Value = "  +6.23%  "
$ws.Range("D12").Value = "'35.25"
$ws.Range("E12").Value = "  +8.83%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'16.16"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "'0.910"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "2.357.51"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "43.783.20"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'77.64"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").Value = "'252.87"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("E25").Value = "  -6.27%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "'175.97"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'22.25"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'0.132"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("D33").Value = "'0.0749"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'6.60"
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  +13.48%  "
$ws.Range("D41").Value = "'20.26"
$ws.Range("E41").Value = "  +6.45%  "
$ws.Range("D42").Value = "'64.44"
$ws.Range("E42").Value = "  +10.68%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.106"
$ws.Range("E44").Value = "  -6.78%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.02"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'2.47"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "'97.72"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("E51").Value = "  +1.96%  "
